$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update odds values that changed for rows 2, 3, 4 and 7

# Row 2
$ws.Range("G2").Value = 1.73
$ws.Range("H2").Value = 3.4
$ws.Range("I2").Value = 5.75
$ws.Range("J2").Value = 2.4
$ws.Range("K2").Value = 2
$ws.Range("Q2").Value = 2.5
$ws.Range("R2").Value = 1.5
$ws.Range("X2").Value = 6.5
$ws.Range("Y2").Value = 9
$ws.Range("AA2").Value = 17
$ws.Range("AG2").Value = 11
$ws.Range("AH2").Value = 26
$ws.Range("AJ2").Value = 67
$ws.Range("AN2").Value = 3.4
$ws.Range("AO2").Value = 9.5
$ws.Range("AQ2").Value = 34
$ws.Range("AS2").Value = 251
$ws.Range("AW2").Value = 7
$ws.Range("AZ2").Value = 151
$ws.Range("BD2").Value = 151

# Row 3
$ws.Range("G3").Value = 2.75
$ws.Range("H3").Value = 3.1
$ws.Range("I3").Value = 2.8
$ws.Range("L3").Value = 3.6
$ws.Range("M3").Value = 1.1
$ws.Range("N3").Value = 7
$ws.Range("O3").Value = 1.5
$ws.Range("P3").Value = 2.5
$ws.Range("Q3").Value = 2.6
$ws.Range("R3").Value = 1.48
$ws.Range("S3").Value = 1.57
$ws.Range("T3").Value = 2.25
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("X3").Value = 12
$ws.Range("AC3").Value = 6.5
$ws.Range("AF3").Value = 67
$ws.Range("AI3").Value = 11
$ws.Range("AK3").Value = 26
$ws.Range("AS3").Value = 301
$ws.Range("AT3").Value = 2.25
$ws.Range("AU3").Value = 9
$ws.Range("AW3").Value = 4.5
$ws.Range("AX3").Value = 17
$ws.Range("AZ3").Value = 51
$ws.Range("BB3").Value = 301

# Row 4
$ws.Range("G4").Value = 3.7
$ws.Range("H4").Value = 3
$ws.Range("K4").Value = 1.95
$ws.Range("O4").Value = 1.44
$ws.Range("P4").Value = 2.63
$ws.Range("AD4").Value = 6
$ws.Range("AQ4").Value = 67
$ws.Range("AR4").Value = 101
$ws.Range("AS4").Value = 301
$ws.Range("AZ4").Value = 41

# Row 7
$ws.Range("M7").Value = 1.1
$ws.Range("N7").Value = 7
$ws.Range("Q7").Value = 2.4
$ws.Range("R7").Value = 1.53

# Remove row 8 (Paraguay Primera Division: Ameliano - 2 de Mayo) entirely
$ws.Rows.Item(8).Delete()

